$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.270858666666666
$ws.Range("H2").Value = 15.812576
$ws.Range("I2").Value = 0.02659672411376159
$ws.Range("J2").Value = 0.02659672411376159
$ws.Range("M2").Value = 44.40220133333333
$ws.Range("N2").Value = 133.206604
$ws.Range("O2").Value = 0.9893265572082102
$ws.Range("P2").Value = 0.9893265572082101
$ws.Range("Q2").Value = 234.0377277168782
$ws.Range("R2").Value = 2106.339549451904
$ws.Range("S2").Value = 0.02631284550048434
$ws.Range("T2").Value = 0.02631284550048434
# Row 3
$ws.Range("G3").Value = 5.270858666666666
$ws.Range("H3").Value = 15.812576
$ws.Range("I3").Value = 0.02659672411376159
$ws.Range("J3").Value = 0.02659672411376159
$ws.Range("M3").Value = 0.401961
$ws.Range("N3").Value = 1.205883
$ws.Range("O3").Value = 0.008956103083191794
$ws.Range("P3").Value = 0.008956103083191792
$ws.Range("Q3").Value = 2.118679620512
$ws.Range("R3").Value = 19.068116584608
$ws.Range("S3").Value = 0.0002382030028380617
$ws.Range("T3").Value = 0.0002382030028380616
# Row 4
$ws.Range("G4").Value = 5.270858666666666
$ws.Range("H4").Value = 15.812576
$ws.Range("I4").Value = 0.02659672411376159
$ws.Range("J4").Value = 0.02659672411376159
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07707633333333333
$ws.Range("N4").Value = 0.231229
$ws.Range("O4").Value = 0.00171733970859806
$ws.Range("P4").Value = 0.00171733970859806
$ws.Range("Q4").Value = 0.4062584595448889
$ws.Range("R4").Value = 3.656326135904
$ws.Range("S4").Value = 0.00004567561043919034
$ws.Range("T4").Value = 0.00004567561043919033
# Row 5
$ws.Range("I5").Value = 0.01269587703542364
$ws.Range("J5").Value = 0.01269587703542364
$ws.Range("M5").Value = 44.40220133333333
$ws.Range("N5").Value = 133.206604
$ws.Range("O5").Value = 0.9893265572082102
$ws.Range("P5").Value = 0.9893265572082101
$ws.Range("Q5").Value = 111.7173002221742
$ws.Range("R5").Value = 1005.455701999568
$ws.Range("S5").Value = 0.01256036831819444
$ws.Range("T5").Value = 0.01256036831819444
# Row 6
$ws.Range("I6").Value = 0.01269587703542364
$ws.Range("J6").Value = 0.01269587703542364
$ws.Range("M6").Value = 0.401961
$ws.Range("N6").Value = 1.205883
$ws.Range("O6").Value = 0.008956103083191794
$ws.Range("P6").Value = 0.008956103083191792
$ws.Range("Q6").Value = 1.011346202804
$ws.Range("R6").Value = 9.102115825236
$ws.Range("S6").Value = 0.0001137055834607815
$ws.Range("T6").Value = 0.0001137055834607815
# Row 7
$ws.Range("I7").Value = 0.01269587703542364
$ws.Range("J7").Value = 0.01269587703542364
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07707633333333333
$ws.Range("N7").Value = 0.231229
$ws.Range("O7").Value = 0.00171733970859806
$ws.Range("P7").Value = 0.00171733970859806
$ws.Range("Q7").Value = 0.1939264183408889
$ws.Range("R7").Value = 1.745337765068
$ws.Range("S7").Value = 0.00002180313376841124
$ws.Range("T7").Value = 0.00002180313376841123
# Row 8
$ws.Range("G8").Value = 114.018682
$ws.Range("H8").Value = 342.056046
$ws.Range("I8").Value = 0.5753376481419691
$ws.Range("J8").Value = 0.5753376481419691
$ws.Range("M8").Value = 44.40220133333333
$ws.Range("N8").Value = 133.206604
$ws.Range("O8").Value = 0.9893265572082102
$ws.Range("P8").Value = 0.9893265572082101
$ws.Range("Q8").Value = 5062.680473925309
$ws.Range("R8").Value = 45564.12426532779
$ws.Range("S8").Value = 0.5691968146685629
$ws.Range("T8").Value = 0.5691968146685628
# Row 9
$ws.Range("G9").Value = 114.018682
$ws.Range("H9").Value = 342.056046
$ws.Range("I9").Value = 0.5753376481419691
$ws.Range("J9").Value = 0.5753376481419691
$ws.Range("M9").Value = 0.401961
$ws.Range("N9").Value = 1.205883
$ws.Range("O9").Value = 0.008956103083191794
$ws.Range("P9").Value = 0.008956103083191792
$ws.Range("Q9").Value = 45.83106343540201
$ws.Range("R9").Value = 412.479570918618
$ws.Range("S9").Value = 0.005152783284400605
$ws.Range("T9").Value = 0.005152783284400604
# Row 10
$ws.Range("G10").Value = 114.018682
$ws.Range("H10").Value = 342.056046
$ws.Range("I10").Value = 0.5753376481419691
$ws.Range("J10").Value = 0.5753376481419691
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.07707633333333333
$ws.Range("N10").Value = 0.231229
$ws.Range("O10").Value = 0.00171733970859806
$ws.Range("P10").Value = 0.00171733970859806
$ws.Range("Q10").Value = 8.788141940059335
$ws.Range("R10").Value = 79.09327746053401
$ws.Range("S10").Value = 0.0009880501890056227
$ws.Range("T10").Value = 0.0009880501890056224
# Row 11
$ws.Range("G11").Value = 1.265015666666667
$ws.Range("H11").Value = 3.795047
$ws.Range("I11").Value = 0.006383262161570549
$ws.Range("J11").Value = 0.006383262161570549
$ws.Range("M11").Value = 44.40220133333333
$ws.Range("N11").Value = 133.206604
$ws.Range("O11").Value = 0.9893265572082102
$ws.Range("P11").Value = 0.9893265572082101
$ws.Range("Q11").Value = 56.16948032115423
$ws.Range("R11").Value = 505.5253228903881
$ws.Range("S11").Value = 0.006315130778064029
$ws.Range("T11").Value = 0.006315130778064029
# Row 12
$ws.Range("G12").Value = 1.265015666666667
$ws.Range("H12").Value = 3.795047
$ws.Range("I12").Value = 0.006383262161570549
$ws.Range("J12").Value = 0.006383262161570549
$ws.Range("M12").Value = 0.401961
$ws.Range("N12").Value = 1.205883
$ws.Range("O12").Value = 0.008956103083191794
$ws.Range("P12").Value = 0.008956103083191792
$ws.Range("Q12").Value = 0.508486962389
$ws.Range("R12").Value = 4.576382661501
$ws.Range("S12").Value = 0.00005716915392606351
$ws.Range("T12").Value = 0.00005716915392606350
# Row 13
$ws.Range("G13").Value = 1.265015666666667
$ws.Range("H13").Value = 3.795047
$ws.Range("I13").Value = 0.006383262161570549
$ws.Range("J13").Value = 0.006383262161570549
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.07707633333333333
$ws.Range("N13").Value = 0.231229
$ws.Range("O13").Value = 0.00171733970859806
$ws.Range("P13").Value = 0.00171733970859806
$ws.Range("Q13").Value = 0.09750276919588889
$ws.Range("R13").Value = 0.877524922763
$ws.Range("S13").Value = 0.00001096222958045659
$ws.Range("T13").Value = 0.00001096222958045659
# Row 14
$ws.Range("G14").Value = 75.10640066666666
$ws.Range("H14").Value = 225.319202
$ws.Range("I14").Value = 0.3789864885472752
$ws.Range("J14").Value = 0.3789864885472752
$ws.Range("M14").Value = 44.40220133333333
$ws.Range("N14").Value = 133.206604
$ws.Range("O14").Value = 0.9893265572082102
$ws.Range("P14").Value = 0.9893265572082101
$ws.Range("Q14").Value = 3334.889523823334
$ws.Range("R14").Value = 30014.00571441001
$ws.Range("S14").Value = 0.3749413979429045
$ws.Range("T14").Value = 0.3749413979429045
# Row 15
$ws.Range("G15").Value = 75.10640066666666
$ws.Range("H15").Value = 225.319202
$ws.Range("I15").Value = 0.3789864885472752
$ws.Range("J15").Value = 0.3789864885472752
$ws.Range("M15").Value = 0.401961
$ws.Range("N15").Value = 1.205883
$ws.Range("O15").Value = 0.008956103083191794
$ws.Range("P15").Value = 0.008956103083191792
$ws.Range("Q15").Value = 30.189843918374
$ws.Range("R15").Value = 271.708595265366
$ws.Range("S15").Value = 0.003394242058566282
$ws.Range("T15").Value = 0.003394242058566282
# Row 16
$ws.Range("G16").Value = 75.10640066666666
$ws.Range("H16").Value = 225.319202
$ws.Range("I16").Value = 0.3789864885472752
$ws.Range("J16").Value = 0.3789864885472752
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07707633333333333
$ws.Range("N16").Value = 0.231229
$ws.Range("O16").Value = 0.00171733970859806
$ws.Range("P16").Value = 0.00171733970859806
$ws.Range("Q16").Value = 8.788141940059335
$ws.Range("R16").Value = 79.09327746053401
$ws.Range("S16").Value = 0.0009880501890056227
$ws.Range("T16").Value = 0.0009880501890056224
